# fix: apply menu type enum to the all menu service codes
#
# The workbook's single sheet is renamed from the week-specific
# "12월3주" ("3rd week of December") to the generic "식단표" ("menu
# table") so the template works for any week/menu code, and the
# worksheet's stored selection is moved from the old E12:E16 block to
# the title block A1:F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "12월3주" -> "식단표"
$ws.Name = "식단표"

# Move/update the active selection to the title merge range A1:F4
$ws.Range("A1:F4").Select()
